$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.86003550082103
$ws.Range("D2").Value = 4.446823480900126
$ws.Range("E2").Value = 19.1350167174431
$ws.Range("F2").Value = 21.59066055625179
$ws.Range("G2").Value = 24.14616097514553
$ws.Range("H2").Value = 13.11496453583104
$ws.Range("K2").Value = 10.48029390352347
$ws.Range("L2").Value = 8.056990899017634
$ws.Range("N2").Value = 19.20809583049475
$ws.Range("O2").Value = 19.347415422169

$ws.Range("B3").Value = 15.74898033238184
$ws.Range("D3").Value = 4.368082632537432
$ws.Range("E3").Value = 19.20872617613973
$ws.Range("F3").Value = 21.58745153567825
$ws.Range("G3").Value = 24.13080695176284
$ws.Range("H3").Value = 13.14999524576426
$ws.Range("K3").Value = 10.22289605073004
$ws.Range("L3").Value = 8.011465306176349
$ws.Range("N3").Value = 19.26568275174242
$ws.Range("O3").Value = 19.39269144727861

$ws.Range("B4").Value = 15.68361088360445
$ws.Range("D4").Value = 4.318304837739656
$ws.Range("E4").Value = 19.257099452887
$ws.Range("F4").Value = 21.59113795121808
$ws.Range("G4").Value = 24.12926143830791
$ws.Range("H4").Value = 13.1734844166666
$ws.Range("K4").Value = 10.05991937796951
$ws.Range("L4").Value = 7.98471019097011
$ws.Range("N4").Value = 19.30278746372706
$ws.Range("O4").Value = 19.42454672836973

$ws.Range("B5").Value = 15.65770598458786
$ws.Range("D5").Value = 4.297673466523229
$ws.Range("E5").Value = 19.27759585435717
$ws.Range("F5").Value = 21.59406388031546
$ws.Range("G5").Value = 24.13061444919467
$ws.Range("H5").Value = 13.18355438506508
$ws.Range("K5").Value = 9.99232299846066
$ws.Range("L5").Value = 7.974117201565631
$ws.Range("N5").Value = 19.31834820735479
$ws.Range("O5").Value = 19.43854630921661

$ws.Range("B6").Value = 15.65344951095414
$ws.Range("D6").Value = 4.294227112343525
$ws.Range("E6").Value = 19.28104662714198
$ws.Range("F6").Value = 21.59463573313251
$ws.Range("G6").Value = 24.13095886146062
$ws.Range("H6").Value = 13.18525657008069
$ws.Range("K6").Value = 9.981029064693026
$ws.Range("L6").Value = 7.972377210242858
$ws.Range("N6").Value = 19.32095868367806
$ws.Range("O6").Value = 19.44093238358594

$ws.Range("B7").Value = 15.68325851848888
$ws.Range("D7").Value = 4.318027980987384
$ws.Range("E7").Value = 19.25737269984706
$ws.Range("F7").Value = 21.59117164594883
$ws.Range("G7").Value = 24.12927165738733
$ws.Range("H7").Value = 13.17361820763761
$ws.Range("K7").Value = 10.05901245549659
$ws.Range("L7").Value = 7.984566064052077
$ws.Range("N7").Value = 19.30299553711764
$ws.Range("O7").Value = 19.42473141071811

$ws.Range("B8").Value = 15.8211738911381
$ws.Range("D8").Value = 4.4199782964257
$ws.Range("E8").Value = 19.15978545159089
$ws.Range("F8").Value = 21.58838105505928
$ws.Range("G8").Value = 24.13923207985958
$ws.Range("H8").Value = 13.126632150082
$ws.Range("K8").Value = 10.39260086854923
$ws.Range("L8").Value = 8.041050448660632
$ws.Range("N8").Value = 19.22759011966096
$ws.Range("O8").Value = 19.36218398585948

$ws.Range("B9").Value = 16.1128461400389
$ws.Range("D9").Value = 4.608003343921716
$ws.Range("E9").Value = 18.99311298859554
$ws.Range("F9").Value = 21.62769698746067
$ws.Range("G9").Value = 24.22118299948126
$ws.Range("H9").Value = 13.05020372960604
$ws.Range("K9").Value = 11.00526144734008
$ws.Range("L9").Value = 8.160917113451983
$ws.Range("N9").Value = 19.09351958477334
$ws.Range("O9").Value = 19.27177022361838

$ws.Range("B10").Value = 16.33847180496899
$ws.Range("D10").Value = 4.738204561387849
$ws.Range("E10").Value = 18.88568171530256
$ws.Range("F10").Value = 21.68369778966754
$ws.Range("G10").Value = 24.31915926734496
$ws.Range("H10").Value = 13.00362829442399
$ws.Range("K10").Value = 11.42720709992302
$ws.Range("L10").Value = 8.253967150870317
$ws.Range("N10").Value = 19.00335200699507
$ws.Range("O10").Value = 19.22507413816395

$ws.Range("B11").Value = 16.44321844566346
$ws.Range("D11").Value = 4.79558006674578
$ws.Range("E11").Value = 18.84006405863349
$ws.Range("F11").Value = 21.71500055518838
$ws.Range("G11").Value = 24.37183069717235
$ws.Range("H11").Value = 12.9845188279126
$ws.Range("K11").Value = 11.61251479184965
$ws.Range("L11").Value = 8.297253547970488
$ws.Range("N11").Value = 18.96412546230664
$ws.Range("O11").Value = 19.20812748161552

$ws.Range("B12").Value = 16.48315559712375
$ws.Range("D12").Value = 4.817029911388857
$ws.Range("E12").Value = 18.82325726833595
$ws.Range("F12").Value = 21.72768523590266
$ws.Range("G12").Value = 24.39292935733266
$ws.Range("H12").Value = 12.9775813201187
$ws.Range("K12").Value = 11.68169243762326
$ws.Range("L12").Value = 8.313771599856363
$ws.Range("N12").Value = 18.94952773073095
$ws.Range("O12").Value = 19.20232855488386

$ws.Range("B13").Value = 16.47454282739688
$ws.Range("D13").Value = 4.812422783131791
$ws.Range("E13").Value = 18.82685611658847
$ws.Range("F13").Value = 21.72491652649316
$ws.Range("G13").Value = 24.38833430650927
$ws.Range("H13").Value = 12.97906214566796
$ws.Range("K13").Value = 11.66683861440995
$ws.Range("L13").Value = 8.31020869979853
$ws.Range("N13").Value = 18.95266022172143
$ws.Range("O13").Value = 19.20354994514372

$ws.Range("B14").Value = 16.446498830944
$ws.Range("D14").Value = 4.79735036763755
$ws.Range("E14").Value = 18.83867198305489
$ws.Range("F14").Value = 21.71602752080095
$ws.Range("G14").Value = 24.37354346082773
$ws.Range("H14").Value = 12.98394208583574
$ws.Range("K14").Value = 11.61822622225499
$ws.Range("L14").Value = 8.298610029239523
$ws.Range("N14").Value = 18.96291936329999
$ws.Range("O14").Value = 19.20763800349711

$ws.Range("B15").Value = 16.42935554411712
$ws.Range("D15").Value = 4.788081703748894
$ws.Range("E15").Value = 18.84597043545111
$ws.Range("F15").Value = 21.71069074350996
$ws.Range("G15").Value = 24.36463342825446
$ws.Range("H15").Value = 12.98697010738485
$ws.Range("K15").Value = 11.58831911261477
$ws.Range("L15").Value = 8.291521636796094
$ws.Range("N15").Value = 18.96923675780977
$ws.Range("O15").Value = 19.21022260858084

$ws.Range("B16").Value = 16.33166612338306
$ws.Range("D16").Value = 4.7344167270091
$ws.Range("E16").Value = 18.88872838465121
$ws.Range("F16").Value = 21.68176877264908
$ws.Range("G16").Value = 24.31587911913804
$ws.Range("H16").Value = 13.00491895920316
$ws.Range("K16").Value = 11.4149597488405
$ws.Range("L16").Value = 8.251156592505778
$ws.Range("N16").Value = 19.00595150582638
$ws.Range("O16").Value = 19.22626814594544

$ws.Range("B17").Value = 16.27225467841059
$ws.Range("D17").Value = 4.701012574336499
$ws.Range("E17").Value = 18.9157921159096
$ws.Range("F17").Value = 21.66551403403419
$ws.Range("G17").Value = 24.28803709898402
$ws.Range("H17").Value = 13.01646223340758
$ws.Range("K17").Value = 11.30687981331581
$ws.Range("L17").Value = 8.22663123530484
$ws.Range("N17").Value = 19.02893281392495
$ws.Range("O17").Value = 19.23721234581959

$ws.Range("B18").Value = 16.23828314042078
$ws.Range("D18").Value = 4.681625493973707
$ws.Range("E18").Value = 18.93166469176862
$ws.Range("F18").Value = 21.65671375647754
$ws.Range("G18").Value = 24.27278649280239
$ws.Range("H18").Value = 13.0232972038335
$ws.Range("K18").Value = 11.24409248490434
$ws.Range("L18").Value = 8.21261570716722
$ws.Range("N18").Value = 19.04231970710902
$ws.Range("O18").Value = 19.24391145737744

$ws.Range("B19").Value = 16.22681632041143
$ws.Range("D19").Value = 4.675031824733434
$ws.Range("E19").Value = 18.93709147564702
$ws.Range("F19").Value = 21.65382862504077
$ws.Range("G19").Value = 24.26775434205113
$ws.Range("H19").Value = 13.02564499799967
$ws.Range("K19").Value = 11.222728116785
$ws.Range("L19").Value = 8.207886226166321
$ws.Range("N19").Value = 19.04688127771093
$ws.Range("O19").Value = 19.24624907487748

$ws.Range("B20").Value = 16.27855861799488
$ws.Range("D20").Value = 4.704586578807925
$ws.Range("E20").Value = 18.91287944304841
$ws.Range("F20").Value = 21.66718760494801
$ws.Range("G20").Value = 24.29092200164941
$ws.Range("H20").Value = 13.01521318963652
$ws.Range("K20").Value = 11.31844984517046
$ws.Range("L20").Value = 8.22923268397293
$ws.Range("N20").Value = 19.02646896720823
$ws.Range("O20").Value = 19.23600546802081

$ws.Range("B21").Value = 16.45472890513285
$ws.Range("D21").Value = 4.80178509948868
$ws.Range("E21").Value = 18.83518868694387
$ws.Range("F21").Value = 21.71861594506365
$ws.Range("G21").Value = 24.37785669898203
$ws.Range("H21").Value = 12.98250061876788
$ws.Range("K21").Value = 11.63253214650686
$ws.Range("L21").Value = 8.302013499169423
$ws.Range("N21").Value = 18.95989905147008
$ws.Range("O21").Value = 19.20642045392167

$ws.Range("B22").Value = 16.57143632839511
$ws.Range("D22").Value = 4.863690927147589
$ws.Range("E22").Value = 18.78713888905288
$ws.Range("F22").Value = 21.75706708444926
$ws.Range("G22").Value = 24.4413880571585
$ws.Range("H22").Value = 12.96286288509074
$ws.Range("K22").Value = 11.83199126315709
$ws.Range("L22").Value = 8.35031157647685
$ws.Range("N22").Value = 18.91788641778581
$ws.Range("O22").Value = 19.19068963081565

$ws.Range("B23").Value = 16.50901415034388
$ws.Range("D23").Value = 4.830802042800523
$ws.Range("E23").Value = 18.81253464020169
$ws.Range("F23").Value = 21.73610473236739
$ws.Range("G23").Value = 24.40687020431818
$ws.Range("H23").Value = 12.97318453233688
$ws.Range("K23").Value = 11.72608020248882
$ws.Range("L23").Value = 8.32447080392026
$ws.Range("N23").Value = 18.94017293482911
$ws.Range("O23").Value = 19.19875546644029

$ws.Range("B24").Value = 16.27570802860712
$ws.Range("D24").Value = 4.702971339309597
$ws.Range("E24").Value = 18.91419528661682
$ws.Range("F24").Value = 21.66642928605334
$ws.Range("G24").Value = 24.28961538096096
$ws.Range("H24").Value = 13.01577726373663
$ws.Range("K24").Value = 11.31322105761922
$ws.Range("L24").Value = 8.228056305202951
$ws.Range("N24").Value = 19.02758232825434
$ws.Range("O24").Value = 19.23654982930063

$ws.Range("B25").Value = 16.03183131241305
$ws.Range("D25").Value = 4.558484934120863
$ws.Range("E25").Value = 19.03556185650563
$ws.Range("F25").Value = 21.61228153020159
$ws.Range("G25").Value = 24.19235250346236
$ws.Range("H25").Value = 13.06919773992295
$ws.Range("K25").Value = 10.84427621173226
$ws.Range("L25").Value = 8.127572581506151
$ws.Range("N25").Value = 19.12832006018177
$ws.Range("O25").Value = 19.29276926642135
